$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price/Volume columns so numeric-looking strings are not
# auto-converted to numbers by Excel (matches original inlineStr text cells).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '60.783.53'
$ws.Range('E2').Value = '  -2.00%  '

$ws.Range('D3').Value = '3.382.05'
$ws.Range('E3').Value = '  -0.97%  '

$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').Value = '574.48'
$ws.Range('E5').Value = '  -0.66%  '

$ws.Range('D6').Value = '136.65'
$ws.Range('E6').Value = '  -1.60%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('D8').Value = '3.381.46'
$ws.Range('E8').Value = '  -0.95%  '

$ws.Range('D9').Value = '0.471'
$ws.Range('E9').Value = '  -1.55%  '

$ws.Range('D10').Value = '7.58'
$ws.Range('E10').Value = '  +1.10%  '

$ws.Range('D11').Value = '0.123'
$ws.Range('E11').Value = '  -3.43%  '

$ws.Range('D12').Value = '0.390'
$ws.Range('E12').Value = '  -1.28%  '

$ws.Range('D13').Value = '3.966.91'
$ws.Range('E13').Value = '  -0.83%  '

$ws.Range('E14').Value = '  +0.91%  '

$ws.Range('D15').Value = '26.38'
$ws.Range('E15').Value = '  +3.42%  '

$ws.Range('D16').Value = '0.0000173'
$ws.Range('E16').Value = '  -4.06%  '

$ws.Range('D17').Value = '3.392.38'
$ws.Range('E17').Value = '  -0.66%  '

$ws.Range('D18').Value = '60.941.74'
$ws.Range('E18').Value = '  -1.71%  '

$ws.Range('D19').Value = '14.03'
$ws.Range('E19').Value = '  -1.20%  '

$ws.Range('D20').Value = '5.84'
$ws.Range('E20').Value = '  -1.04%  '

$ws.Range('D21').Value = '9.44'
$ws.Range('E21').Value = '  -1.15%  '

$ws.Range('D22').Value = '377.36'
$ws.Range('E22').Value = '  -3.28%  '

$ws.Range('D23').Value = '0.557'
$ws.Range('E23').Value = '  -2.85%  '

$ws.Range('D24').Value = '3.529.11'
$ws.Range('E24').Value = '  -0.62%  '

$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.15%  '

$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '71.40'
$ws.Range('E26').Value = '  -0.26%  '

$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '0.0000125'
$ws.Range('E27').Value = '  -2.69%  '

$ws.Range('D28').Value = '1.79'
$ws.Range('E28').Value = '  +12.52%  '

$ws.Range('D29').Value = '7.57'
$ws.Range('E29').Value = '  -1.38%  '

$ws.Range('E30').Value = '  +4.99%  '

$ws.Range('E31').Value = '  +0.04%  '

$ws.Range('D32').Value = '8.19'
$ws.Range('E32').Value = '  -1.47%  '

$ws.Range('D33').Value = '2.16'
$ws.Range('E33').Value = '  -0.79%  '

$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('D35').Value = '23.77'
$ws.Range('E35').Value = '  +0.71%  '

$ws.Range('D36').Value = '5.21'
$ws.Range('E36').Value = '  -5.18%  '

$ws.Range('D37').Value = '6.87'
$ws.Range('E37').Value = '  -2.11%  '

$ws.Range('D38').Value = '1.54'
$ws.Range('E38').Value = '  -1.80%  '

$ws.Range('D39').Value = '164.58'
$ws.Range('E39').Value = '  +1.03%  '

$ws.Range('D40').Value = '0.0758'
$ws.Range('E40').Value = '  -4.34%  '

$ws.Range('E41').Value = '  +0.06%  '

$ws.Range('D42').Value = '0.772'
$ws.Range('E42').Value = '  -2.37%  '

$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '1.71'
$ws.Range('E43').Value = '  -3.21%  '

$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = '4.42'
$ws.Range('E44').Value = '  -1.39%  '

$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '41.53'
$ws.Range('E45').Value = '  -0.52%  '

$ws.Range('B46').Value = 'ONDO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D46').Value = '1.20'
$ws.Range('E46').Value = '  -2.85%  '

$ws.Range('D47').Value = '24.47'
$ws.Range('E47').Value = '  -2.60%  '

$ws.Range('D48').Value = '23.45'
$ws.Range('E48').Value = '  +1.37%  '

$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = '6.80'
$ws.Range('E49').Value = '  -2.58%  '

$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.409.41'
$ws.Range('E50').Value = '  +1.37%  '

$ws.Range('D51').Value = '2.42'
$ws.Range('E51').Value = '  +5.67%  '
